$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.101.67'
$ws.Range('E2').Value = '  +2.44%  '
$ws.Range('D3').Value = '3.592.55'
$ws.Range('E3').Value = '  +1.23%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.35%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '586.25'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +2.70%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '186.51'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +1.43%  '
$ws.Range('D7').Value = '3.581.36'
$ws.Range('E7').Value = '  +1.17%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.622'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +1.16%  '
$ws.Range('E9').Value = '  -0.03%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.218'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +19.40%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.652'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +1.38%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '54.49'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +1.36%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.0000321'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +7.56%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '9.56'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +0.85%  '
$ws.Range('D15').Value = '4.159.59'
$ws.Range('E15').Value = '  +0.96%  '
$ws.Range('D16').Value = '70.947.84'
$ws.Range('E16').Value = '  +2.22%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '19.30'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').Value = '3.585.95'
$ws.Range('E18').Value = '  +0.52%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '571.13'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +14.33%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '12.41'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +0.73%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.121'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('E22').Value = '  -2.08%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '17.55'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -9.61%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '5.06'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +3.59%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '4.62'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +7.35%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '95.14'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +1.31%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '11.31'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.17%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.94'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +0.46%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '9.13'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -0.73%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '32.40'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +3.23%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '7.24'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -4.13%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '12.31'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.36%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.115'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.00%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '64.20'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -1.84%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '3.32'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +6.32%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '548.93'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -2.92%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.414'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +3.12%  '
$ws.Range('D38').Value = '0.0₃0814'
$ws.Range('E38').Value = '  +4.58%  '
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '37.67'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.44%  '
$ws.Range('D41').Value = '3.494.61'
$ws.Range('E41').Value = '  +10.15%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '3.20'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -3.15%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '3.46'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +2.37%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.136'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +1.86%  '
$ws.Range('E45').Value = '  -1.61%  '
$ws.Range('E46').Value = '  -0.36%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0444'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +0.47%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '9.42'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +1.53%  '
$ws.Range('E49').Value = '  +2.77%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.997'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.54%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '1.44'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -0.93%  '
